# Project Admin flow update: duplicate the "Travel Honoraria Terms" (THT_*)
# 4-column block (BU:BX) into four new trailing columns (BY:CB), fix up the
# BU2 data value, and populate the new row-2 cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): clone BU1:BX1 headers into BY1:CB1 ---------------
$ws.Range("BY1").Value = $ws.Range("BU1").Value2
$ws.Range("BZ1").Value = $ws.Range("BV1").Value2
$ws.Range("CA1").Value = $ws.Range("BW1").Value2
$ws.Range("CB1").Value = $ws.Range("BX1").Value2

# --- Data row (row 2) -------------------------------------------------------
# BU2 changes from 2 -> 1
$ws.Range("BU2").Value = 1

# New trailing columns, all populated with 1
$ws.Range("BY2").Value = 1
$ws.Range("BZ2").Value = 1
$ws.Range("CA2").Value = 1
$ws.Range("CB2").Value = 1

# --- View state: move the viewport / selection to match the new layout ----
$win = $excel.ActiveWindow
[void]$ws.Range("BW2").Select()
$win.ScrollColumn = 61
$win.ScrollRow = 1
